$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "nat_escapement" column (B) is being removed entirely - its data will
# be sourced from a different file. Remove the threaded comment that was
# attached to its header cell (B1) before shifting columns, then delete the
# whole column so ocean_catch/ocean_er shift left from C/D into B/C.
$commentB1 = $ws.Range("B1").Comment
if ($commentB1 -ne $null) {
    $commentB1.Delete()
}

$ws.Columns("B").Delete()

# Restore a sensible active selection near where the data now ends.
$ws.Range("I46").Select()
